# Applies the "changed MP time limit and corrected error in fixed recourse data" edit.
#
# Sheet1 (summary, rows 2-11, one row per instance 1..10): update objective (B),
# solve time (C), num cuts (F), num variables (G), num cons (H), num quad cons (I).
#
# Per-instance tabs "1".."10" (iteration log, rows 2-3): update MP solve time (D2),
# Worst violation (E2), MP objective (B3), MP gap (C3), MP solve time (D3).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: summary table -------------------------------------------------
$summary = @{
    2  = @{ B = -105.56745778630673; C = 11.123585724; F = 20; G = 4540; H = 4900; I = 400 }
    3  = @{ B = -105.84488399353643; C = 5.559538512;  F = 20; G = 4540; H = 4900; I = 400 }
    4  = @{ B = -105.49382704393152; C = 1.652208195;  F = 20; G = 4540; H = 4900; I = 400 }
    5  = @{ B = -104.60548318016075; C = 1.985973199;  F = 20; G = 4540; H = 4900; I = 400 }
    6  = @{ B = -103.49903626996942; C = 2.483046185;  F = 20; G = 4540; H = 4900; I = 400 }
    7  = @{ B = -104.63902673475309; C = 2.135797431;  F = 20; G = 4540; H = 4900; I = 400 }
    8  = @{ B = -102.87894890275867; C = 1.806903095;  F = 20; G = 4540; H = 4900; I = 400 }
    9  = @{ B = -105.10124780728994; C = 1.970574475;  F = 20; G = 4540; H = 4900; I = 400 }
    10 = @{ B = -105.99608264981093; C = 1.964038455;  F = 20; G = 4540; H = 4900; I = 400 }
    11 = @{ B = -102.78699673530132; C = 1.340570307;  F = 20; G = 4540; H = 4900; I = 400 }
}

foreach ($row in $summary.Keys) {
    $vals = $summary[$row]
    $ws1.Range("B$row").Value = $vals.B
    $ws1.Range("C$row").Value = $vals.C
    $ws1.Range("F$row").Value = $vals.F
    $ws1.Range("G$row").Value = $vals.G
    $ws1.Range("H$row").Value = $vals.H
    $ws1.Range("I$row").Value = $vals.I
}

# --- Per-instance tabs: iteration log ---------------------------------------
$perTab = @{
    "1"  = @{ D2 = 0.9597090733625488;   E2 = 110.35232; B3 = -105.56745778630673; C3 = 0.0;                  D3 = 1.4601812174804687 }
    "2"  = @{ D2 = 0.019844317255493165; E2 = 111.3416;  B3 = -105.84488399353643; C3 = 0.00948959094675425;  D3 = 5.3511624799942625 }
    "3"  = @{ D2 = 0.011187230572509766; E2 = 111.03335; B3 = -105.49382704393152; C3 = 0.07448691419469279;  D3 = 1.4870565705152587 }
    "4"  = @{ D2 = 0.012217883790039063; E2 = 109.81948; B3 = -104.60548318016075; C3 = 0.0;                  D3 = 1.7787166010548097 }
    "5"  = @{ D2 = 0.025888075358276367; E2 = 113.03886; B3 = -103.49903626996942; C3 = 0.0007685320234566919; D3 = 2.272286114399658 }
    "6"  = @{ D2 = 0.009927427041259766; E2 = 118.31719; B3 = -104.63902673475309; C3 = 0.05628629065809936;  D3 = 1.9807078530994873 }
    "7"  = @{ D2 = 0.011874016303710938; E2 = 108.71757; B3 = -102.87894890275867; C3 = 0.003354180648938238; D3 = 1.6260208998426513 }
    "8"  = @{ D2 = 0.010063803397705079; E2 = 106.85473; B3 = -105.10124780728994; C3 = 0.006043950705532129; D3 = 1.8311555493740235 }
    "9"  = @{ D2 = 0.015876360112426758; E2 = 114.1732;  B3 = -105.99608264981093; D3 = 1.8000432440806884 }
    "10" = @{ D2 = 0.011947974393310547; E2 = 109.25818; B3 = -102.78699673530132; D3 = 1.177843801298462 }
}

foreach ($tabName in $perTab.Keys) {
    # $tabName comes out of the hashtable's .Keys collection as an Int64 even
    # though it was declared as a quoted string literal; force it back to a
    # String so Worksheets.Item() looks the sheet up by name, not by index.
    $tabNameStr = [string]$tabName
    $ws = $wb.Worksheets.Item($tabNameStr)
    $vals = $perTab[$tabName]

    $ws.Range("D2").Value = $vals.D2
    $ws.Range("E2").Value = $vals.E2
    $ws.Range("B3").Value = $vals.B3
    if ($vals.ContainsKey("C3")) {
        $ws.Range("C3").Value = $vals.C3
    }
    $ws.Range("D3").Value = $vals.D3
}
